# Adds an "Email" column (E) for the alumni with a mailto hyperlink,
# per commit "untuk menambahkan input email alumni".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create the mailto hyperlink on E2 first (this is what introduces the
# "naafiridho0505@gmail.com" shared string ahead of the "Email" header
# string, matching the saved workbook's string order), then make sure
# the cell shows just the e-mail address rather than the raw mailto: URL.
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:naafiridho0505@gmail.com")
$ws.Range("E2").Value = "naafiridho0505@gmail.com"

# Header for the new column.
$ws.Range("E1").Value = "Email"

# Match the new column's width.
$ws.Range("E1").EntireColumn.ColumnWidth = 23.3

# Leave the selection on C3, matching the saved workbook state.
$ws.Range("C3").Select()
